$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office theme color values (standard RGB hex -> decimal as R + G*256 + B*65536 for COM RGB())
# dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6 accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000 accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

function ToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$colors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")
for ($i=1; $i -le 12; $i++) {
    $c = $tcs.Colors($i)
    $c.RGB = ToComRGB($colors[$i-1])
}
